$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date in A1 (serial date 45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update prices (doubled)
$ws.Range("D29").Value = 112.4
$ws.Range("D30").Value = 187.2
